# Mise à jour site
# Add a "disponible" (availability) column F with header + "1" markers
# for every product row that has an article/path set (A column not empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: bold, centered text "disponible"
$ws.Range("F1").Value = "disponible"
$ws.Range("F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F1").VerticalAlignment = -4160     # xlTop
$ws.Range("F1").Font.Bold = $true

# Column F width stays the same, but let's make sure column width is set explicitly
$ws.Columns("F").ColumnWidth = 22.77734375

# Rows 2..14: mark availability for rows that already have a value in column A
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop

    $aValue = $ws.Cells.Item($r, 1).Value
    if ($aValue -ne $null -and $aValue -ne "") {
        $cell.Value = "1"
    }
}

# Update the selection to match the authored state
$ws.Range("F8").Select()

$wb.Worksheets(1).Range("A1").Select()
$ws.Range("F8").Select()
